$d = $word.ActiveDocument

# The caption run "In her library, Goessel shows little Bill Hart Jr.'s
# cowboy boots. Behind her at right is the headdress worn by " gets
# "Goessel" expanded to "Dr. Tracey Goessel" (formatting is unchanged
# throughout, so this is a pure text insertion).
$d.Content.Find.Execute(
    "In her library, Goessel shows",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "In her library, Dr. Tracey Goessel shows",
    2
)
